$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1581.7778
$ws.Range("J6").Value = 690
$ws.Range("L6").Value = 2070
$ws.Range("N6").Value = -2294
$ws.Range("H19").Value = 1286
$ws.Range("J19").Value = 1341.5
$ws.Range("L19").Value = 1341.5
$ws.Range("N19").Value = -1691.5
$ws.Range("H40").Value = 3947.6191
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H64").Value = 6971.0386
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 7009.88
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 7009.88
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -7505.88
$ws.Range("H67").Value = 6971.0386
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 7009.88
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 7009.88
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -8725.880000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 181.66667
$ws.Range("I4").Value = 198
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 198
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -82
$ws.Range("N4").Value = -332
$ws.Range("H23").Value = 10964.833
$ws.Range("J23").Value = 10964.833
$ws.Range("L23").Value = 10964.833
$ws.Range("N23").Value = -11482.833
$ws.Range("H37").Value = 24034
$ws.Range("I37").Value = 24034
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 24034
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -23761
$ws.Range("N37").ClearContents()
$ws.Range("H44").Value = 69994.336
$ws.Range("J44").Value = 69994.336
$ws.Range("L44").Value = 69994.336
$ws.Range("N44").Value = -70970.336
$ws.Range("H55").Value = 36347.668
$ws.Range("I55").Value = 19524
$ws.Range("K55").Value = 19524
$ws.Range("M55").Value = -19209
$ws.Range("H63").Value = 3699
$ws.Range("I63").Value = 2499
$ws.Range("K63").Value = 2499
$ws.Range("M63").Value = -1813
$ws.Range("H66").Value = 3699
$ws.Range("I66").Value = 2499
$ws.Range("K66").Value = 12495
$ws.Range("M66").Value = -9063
$ws.Range("H74").Value = 2967.2856
$ws.Range("I74").Value = 2993.611
$ws.Range("K74").Value = 2993.611
$ws.Range("M74").Value = -2119.611
$ws.Range("H77").Value = 2967.2856
$ws.Range("I77").Value = 2993.611
$ws.Range("K77").Value = 14968.055
$ws.Range("M77").Value = -10600.055
$ws.Range("H80").Value = 56250
$ws.Range("I80").Value = 56250
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 56250
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -55252
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 56250
$ws.Range("I83").Value = 56250
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 168750
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -163758
$ws.Range("N83").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 462
$ws.Range("I22").Value = 209.83333
$ws.Range("J22").Value = 966.3333
$ws.Range("K22").Value = 209.83333
$ws.Range("L22").Value = 966.3333
$ws.Range("M22").Value = -36.83332999999999
$ws.Range("N22").Value = -1312.3333
$ws.Range("H51").Value = 59958
$ws.Range("J51").Value = 59958
$ws.Range("L51").Value = 59958
$ws.Range("N51").Value = -60940
$ws.Range("H99").Value = 4108.4375
$ws.Range("I99").Value = 2706.111
$ws.Range("J99").Value = 5911.4287
$ws.Range("K99").Value = 2706.111
$ws.Range("L99").Value = 5911.4287
$ws.Range("M99").Value = -1208.111
$ws.Range("N99").Value = -8907.4287
$ws.Range("H105").Value = 2548.077
$ws.Range("J105").Value = 3385
$ws.Range("L105").Value = 3385
$ws.Range("N105").Value = -6879

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5223.657
$ws.Range("I31").Value = 1780.2222
$ws.Range("J31").Value = 6415.615
$ws.Range("K31").Value = 1780.2222
$ws.Range("L31").Value = 6415.615
$ws.Range("M31").Value = -1485.2222
$ws.Range("N31").Value = -7005.615
$ws.Range("H34").Value = 5223.657
$ws.Range("I34").Value = 1780.2222
$ws.Range("J34").Value = 6415.615
$ws.Range("K34").Value = 1780.2222
$ws.Range("L34").Value = 6415.615
$ws.Range("M34").Value = -1578.2222
$ws.Range("N34").Value = -6819.615
$ws.Range("H58").Value = 3005.4783
$ws.Range("I58").Value = 2878.2563
$ws.Range("K58").Value = 2878.2563
$ws.Range("M58").Value = -2675.2563
$ws.Range("H112").Value = 90996.664
$ws.Range("J112").Value = 90996.664
$ws.Range("L112").Value = 90996.664
$ws.Range("N112").Value = -93950.664
$ws.Range("H122").Value = 5758.933
$ws.Range("I122").Value = 4391.6665
$ws.Range("K122").Value = 13174.9995
$ws.Range("M122").Value = -10724.9995
$ws.Range("H136").Value = 3005.4783
$ws.Range("I136").Value = 2878.2563
$ws.Range("K136").Value = 8634.768899999999
$ws.Range("M136").Value = -6084.768899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 4769039
$ws.Range("I121").Value = 686.6667
$ws.Range("K121").Value = 2060.0001
$ws.Range("M121").Value = -750.0001000000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2253.5334
$ws.Range("I132").Value = 1963.36
$ws.Range("J132").Value = 3704.4
$ws.Range("K132").Value = 5890.08
$ws.Range("L132").Value = 11113.2
$ws.Range("M132").Value = -3360.08
$ws.Range("N132").Value = -16173.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5166.6665
$ws.Range("I68").Value = 6250
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 6250
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -5501
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 5166.6665
$ws.Range("I71").Value = 6250
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 31250
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -27506
$ws.Range("N71").Value = -22488
$ws.Range("H122").Value = 18828.334
$ws.Range("I122").Value = 16455.691
$ws.Range("J122").Value = 24997.2
$ws.Range("K122").Value = 49367.073
$ws.Range("L122").Value = 74991.60000000001
$ws.Range("M122").Value = -46917.073
$ws.Range("N122").Value = -79891.60000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 30026
$ws.Range("I34").Value = 30026
$ws.Range("K34").Value = 30026
$ws.Range("M34").Value = -29823
$ws.Range("H75").Value = 109994.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 109994.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 109994.5
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -111866.5
$ws.Range("H78").Value = 109994.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 109994.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 329983.5
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -339343.5
$ws.Range("H81").Value = 3162.1052
$ws.Range("I81").Value = 2545.4443
$ws.Range("J81").Value = 3717.1
$ws.Range("K81").Value = 5090.8886
$ws.Range("L81").Value = 7434.2
$ws.Range("M81").Value = -4029.8886
$ws.Range("N81").Value = -9556.200000000001
$ws.Range("H84").Value = 3162.1052
$ws.Range("I84").Value = 2545.4443
$ws.Range("J84").Value = 3717.1
$ws.Range("K84").Value = 25454.443
$ws.Range("L84").Value = 37171
$ws.Range("M84").Value = -20150.443
$ws.Range("N84").Value = -47779
$ws.Range("H112").Value = 85499
$ws.Range("J112").Value = 85499
$ws.Range("L112").Value = 85499
$ws.Range("N112").Value = -88453
$ws.Range("H136").Value = 8041390.5
$ws.Range("I136").Value = 889.6667
$ws.Range("K136").Value = 2669.0001
$ws.Range("M136").Value = -119.0001000000002

Write-Host "Applied all Adamantoise_Profits updates"
